$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.701.02'
$ws.Range('E2').Value = '  +2.01%  '
$ws.Range('D3').Value = '1.809.37'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.93'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.75'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +5.41%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.288'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0711'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +7.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0928'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.27%  '
$ws.Range('D12').Value = '2.070.98'
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.05'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -3.00%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.801.29'
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.648'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.95%  '
$ws.Range('D16').Value = '34.727.08'
$ws.Range('E16').Value = '  +2.00%  '
$ws.Range('E17').Value = '  +3.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.76'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '254.74'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.91%  '
$ws.Range('D20').Value = '0.0₃0803'
$ws.Range('E20').Value = '  +8.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.998'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('E22').Value = '  +4.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.24'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.16'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '161.26'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.53'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.18'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.42%  '
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0538'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +4.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.81'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.68'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.35%  '
$ws.Range('E34').Value = '  +3.33%  '
$ws.Range('D35').Value = '1.443.51'
$ws.Range('E35').Value = '  -0.54%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.07'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0193'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.644'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '85.28'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.967'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +7.44%  '
$ws.Range('E41').Value = '  -1.23%  '
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.15'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.11'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +6.97%  '
$ws.Range('E45').Value = '  -0.94%  '
$ws.Range('D46').Value = '1.965.73'
$ws.Range('E46').Value = '  +0.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0491'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '106.03'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +8.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '12.14'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.15%  '
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('D51').Value = '0.0₆0125'
$ws.Range('E51').Value = '  +7.97%  '
